$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their text type (avoid Excel auto-converting
# numeric-looking strings like "580.96" into real numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.269.13'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.089.61'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.96'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.51'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.082.75'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.74%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.528'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.61'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.455'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000245'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.66'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.10%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.30%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.596.76'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.165.42'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.10'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.090.84'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '460.21'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.724'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.43'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.95'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.18'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.17%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.16%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.98'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +8.64%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.21'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.51%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.83'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.61'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0842'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.30'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.01%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.22%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.14'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '433.80'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.72'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.52%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0368'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.863.49'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.58%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.86%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.25'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.98'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.05'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.72%  '
